$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-26 20:49:56"

for ($r = 2; $r -le 411; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
